# Adds "error31r_ref1_permil" (P) and "error31r_ref2_permil" (Q) columns
# to each worksheet, matching the check 31R script error calculation.

$wb = $excel.ActiveWorkbook

# Per-sheet row => error value (permil), same value used for both ref1 and ref2 columns
$dataBySheet = @(
    @{ 2 = -0.4091049258002633; 3 = -0.7444554516042867; 4 = -0.3258536615019514; 5 = -0.661254864741001; 6 = -0.2304847950564959; 7 = -0.2542075049479786; 8 = 0.3000358365834455; 9 = -0.03250962188272588; 10 = 0.009129185265965845; 11 = -0.3232611217129921; 12 = -0.2421222941543233; 13 = 0.804003732371994; 14 = -0.06757790660705254; 15 = 0.7714581516891084; 16 = 0.68887662279149; 17 = 0.4788328094882388; 18 = 0.3962898487157585; 19 = 0.04794438160171666; 20 = 0.07145489776094749; 21 = 1.100060574029005; 22 = 1.123618980839014; 23 = 0.2234770376519091; 24 = 0.2469943301994793 },
    @{ 2 = -0.4091049258002633; 3 = -0.7444554516042867; 4 = -0.3258536615019514; 5 = -0.661254864741001; 6 = -0.2304847950564959; 7 = -0.2542075049479786 },
    @{ 2 = 0.3000358365834455; 3 = -0.03250962188272588; 4 = 0.009129185265965845; 5 = -0.3232611217129921; 6 = -0.2421222941543233; 7 = 0.804003732371994; 8 = -0.06757790660705254 },
    @{ 2 = 0.7714581516891084; 3 = 0.68887662279149; 4 = 0.4788328094882388; 5 = 0.3962898487157585; 6 = 0.04794438160171666; 7 = 0.07145489776094749; 8 = 1.100060574029005; 9 = 1.123618980839014; 10 = 0.2234770376519091; 11 = 0.2469943301994793 }
)

for ($i = 0; $i -lt $dataBySheet.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $vals = $dataBySheet[$i]

    # Header row: copy style from the existing "kappa" header (O1) onto the new headers
    $ws.Range("O1").Copy($ws.Range("P1"))
    $ws.Range("P1").Value = "error31r_ref1_permil"
    $ws.Range("O1").Copy($ws.Range("Q1"))
    $ws.Range("Q1").Value = "error31r_ref2_permil"

    foreach ($row in $vals.Keys) {
        $ws.Range("P$row").Value = $vals[$row]
        $ws.Range("Q$row").Value = $vals[$row]
    }
}
